# Move regression estimates to country specific folder
# Rename the "EL" worksheet to "Retirement_age"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EL")
$ws.Name = "Retirement_age"
